# Apply updates to column A (Regulator Name) on Sheet1 to re-sort/clean
# the translator entries, per the commit "clean and sort out translators
# (sif part done)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = "FOXP3"
    11 = "STAT5"
    17 = "SMAD3"
    18 = "MTOR_DD"
    19 = "NFAT"
    23 = "NFKAPPAB"
    24 = "IL2"
    28 = "CD122"
    29 = "CD132"
    30 = "CD25"
    31 = "IL2R"
    32 = "IL2_EX"
    49 = "PKCTHETA"
    50 = "AKT"
    53 = "TCR"
    54 = "PI3K"
    59 = "CD28"
    60 = "TCR"
    61 = "MTORC2"
    68 = "TCR"
    69 = "IL2_EX"
    70 = "IL2R"
    71 = "RAS"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
